$wb = $excel.ActiveWorkbook
$daily = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

# --- Sheet "Daily": row 2 updates ---
$daily.Range("A2").Value = 47.2229
$daily.Range("B2").Value = 24.7244
$daily.Range("G2").Value = 3538.57
$daily.Range("H2").Value = 6864.02
$daily.Range("I2").Value = 813.05
$daily.Range("J2").Value = 884.64
$daily.Range("L2").Value = 884.64

# --- Sheet "Hourly": per-row updates ---
# row 2
$hourly.Range("A2").Value = 47.2229
$hourly.Range("B2").Value = 24.7244
# row 3
$hourly.Range("A3").Value = 47.2229
$hourly.Range("B3").Value = 24.7244
# row 4
$hourly.Range("A4").Value = 47.2229
$hourly.Range("B4").Value = 24.7244
# row 5
$hourly.Range("A5").Value = 47.2229
$hourly.Range("B5").Value = 24.7244
# row 6
$hourly.Range("A6").Value = 47.2229
$hourly.Range("B6").Value = 24.7244
# row 7
$hourly.Range("A7").Value = 47.2229
$hourly.Range("B7").Value = 24.7244
# row 8
$hourly.Range("A8").Value = 47.2229
$hourly.Range("B8").Value = 24.7244
# row 9
$hourly.Range("A9").Value = 47.2229
$hourly.Range("B9").Value = 24.7244
$hourly.Range("H9").Value = 27.61
$hourly.Range("I9").Value = 163.93
$hourly.Range("J9").Value = 20.93
$hourly.Range("K9").Value = 6.9
$hourly.Range("M9").Value = 6.9
# row 10
$hourly.Range("A10").Value = 47.2229
$hourly.Range("B10").Value = 24.7244
$hourly.Range("H10").Value = 161.89
$hourly.Range("I10").Value = 527.52
$hourly.Range("J10").Value = 59.17
$hourly.Range("K10").Value = 40.47
$hourly.Range("M10").Value = 40.47
# row 11
$hourly.Range("A11").Value = 47.2229
$hourly.Range("B11").Value = 24.7244
$hourly.Range("H11").Value = 311.3
$hourly.Range("I11").Value = 687.63
$hourly.Range("J11").Value = 79.39
$hourly.Range("K11").Value = 77.83
$hourly.Range("M11").Value = 77.83
# row 12
$hourly.Range("A12").Value = 47.2229
$hourly.Range("B12").Value = 24.7244
$hourly.Range("H12").Value = 435.43
$hourly.Range("I12").Value = 769.45
$hourly.Range("J12").Value = 91.58
$hourly.Range("K12").Value = 108.86
$hourly.Range("M12").Value = 108.86
# row 13
$hourly.Range("A13").Value = 47.2229
$hourly.Range("B13").Value = 24.7244
$hourly.Range("H13").Value = 517.6799999999999
$hourly.Range("I13").Value = 811.1
$hourly.Range("J13").Value = 98.40000000000001
$hourly.Range("K13").Value = 129.42
$hourly.Range("M13").Value = 129.42
# row 14
$hourly.Range("A14").Value = 47.2229
$hourly.Range("B14").Value = 24.7244
$hourly.Range("H14").Value = 548.9
$hourly.Range("I14").Value = 825.12
$hourly.Range("J14").Value = 100.77
$hourly.Range("K14").Value = 137.22
$hourly.Range("M14").Value = 137.22
# row 15
$hourly.Range("A15").Value = 47.2229
$hourly.Range("B15").Value = 24.7244
$hourly.Range("H15").Value = 525.8200000000001
$hourly.Range("I15").Value = 814.95
$hourly.Range("J15").Value = 99
$hourly.Range("K15").Value = 131.45
$hourly.Range("M15").Value = 131.45
# row 16
$hourly.Range("A16").Value = 47.2229
$hourly.Range("B16").Value = 24.7244
$hourly.Range("H16").Value = 450.84
$hourly.Range("I16").Value = 778.0700000000001
$hourly.Range("J16").Value = 92.87
$hourly.Range("K16").Value = 112.71
$hourly.Range("M16").Value = 112.71
# row 17
$hourly.Range("A17").Value = 47.2229
$hourly.Range("B17").Value = 24.7244
$hourly.Range("H17").Value = 332.22
$hourly.Range("I17").Value = 703.85
$hourly.Range("J17").Value = 81.59999999999999
$hourly.Range("K17").Value = 83.06
$hourly.Range("M17").Value = 83.06
# row 18
$hourly.Range("A18").Value = 47.2229
$hourly.Range("B18").Value = 24.7244
$hourly.Range("H18").Value = 185.18
$hourly.Range("I18").Value = 560.59
$hourly.Range("J18").Value = 62.93
$hourly.Range("K18").Value = 46.3
$hourly.Range("M18").Value = 46.3
# row 19
$hourly.Range("A19").Value = 47.2229
$hourly.Range("B19").Value = 24.7244
$hourly.Range("H19").Value = 41.69
$hourly.Range("I19").Value = 221.82
$hourly.Range("J19").Value = 26.41
$hourly.Range("K19").Value = 10.42
$hourly.Range("M19").Value = 10.42
# row 20
$hourly.Range("A20").Value = 47.2229
$hourly.Range("B20").Value = 24.7244
# row 21
$hourly.Range("A21").Value = 47.2229
$hourly.Range("B21").Value = 24.7244
# row 22
$hourly.Range("A22").Value = 47.2229
$hourly.Range("B22").Value = 24.7244
# row 23
$hourly.Range("A23").Value = 47.2229
$hourly.Range("B23").Value = 24.7244
# row 24
$hourly.Range("A24").Value = 47.2229
$hourly.Range("B24").Value = 24.7244
# row 25
$hourly.Range("A25").Value = 47.2229
$hourly.Range("B25").Value = 24.7244

# --- sunrise/sunset shared-string text update across all rows ---
$daily.Range("E2").Value = "2024-02-25T07:09:33"
$daily.Range("F2").Value = "2024-02-25T17:59:36"
for ($r = 2; $r -le 25; $r++) {
    $hourly.Cells.Item($r, 5).Value = "2024-02-25T07:09:33"
    $hourly.Cells.Item($r, 6).Value = "2024-02-25T17:59:36"
}
